$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 80

# Copy formatting (styles) from the row above, so the new row keeps the
# same look (bold/border on column A, date format on column E, etc.)
$ws.Range("A79:V79").Copy()
$ws.Range("A80:V80").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 79
$ws.Cells.Item($row, 2).Value = "armenia"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45242.625
$ws.Cells.Item($row, 6).Value = "Ararat Yerevan"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Alashkert"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 4.16
$ws.Cells.Item($row, 11).Value = "11/11/2023 03:12"
$ws.Cells.Item($row, 12).Value = 4.71
$ws.Cells.Item($row, 13).Value = "12/11/2023 14:56"
$ws.Cells.Item($row, 14).Value = 3.41
$ws.Cells.Item($row, 15).Value = "11/11/2023 03:12"
$ws.Cells.Item($row, 16).Value = 3.51
$ws.Cells.Item($row, 17).Value = "12/11/2023 14:56"
$ws.Cells.Item($row, 18).Value = 1.79
$ws.Cells.Item($row, 19).Value = "11/11/2023 03:12"
$ws.Cells.Item($row, 20).Value = 1.78
$ws.Cells.Item($row, 21).Value = "12/11/2023 14:54"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/ararat-yerevan-alashkert/pQNNeYJu/"
